$d = $word.ActiveDocument

# Update the date heading paragraph
$d.Content.Find.Execute("2026-02-22 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2026-02-23 Monday", 2) | Out-Null

# Update the 25 practice-problem table cells (row, col) -> new value
$tbl = $d.Tables.Item(1)

# Cell(1,1): "81÷7=11, 4" -> "76÷3=25, 1"
$tbl.Cell(1, 1).Range.Text = "76÷3=25, 1"

# Cell(1,2): "24÷4=6, 0" -> "76÷6=12, 4"
$tbl.Cell(1, 2).Range.Text = "76÷6=12, 4"

# Cell(1,3): "94÷9=10, 4" -> "15÷3=5, 0"
$tbl.Cell(1, 3).Range.Text = "15÷3=5, 0"

# Cell(1,4): "94÷4=23, 2" -> "36÷7=5, 1"
$tbl.Cell(1, 4).Range.Text = "36÷7=5, 1"

# Cell(1,5): "91÷6=15, 1" -> "59÷9=6, 5"
$tbl.Cell(1, 5).Range.Text = "59÷9=6, 5"

# Cell(5,1): "21÷2=10, 1" -> "86÷5=17, 1"
$tbl.Cell(5, 1).Range.Text = "86÷5=17, 1"

# Cell(5,2): "16÷6=2, 4" -> "38÷7=5, 3"
$tbl.Cell(5, 2).Range.Text = "38÷7=5, 3"

# Cell(5,3): "77÷3=25, 2" -> "55÷2=27, 1"
$tbl.Cell(5, 3).Range.Text = "55÷2=27, 1"

# Cell(5,4): "50÷4=12, 2" -> "94÷3=31, 1"
$tbl.Cell(5, 4).Range.Text = "94÷3=31, 1"

# Cell(5,5): "52÷9=5, 7" -> "82÷7=11, 5"
$tbl.Cell(5, 5).Range.Text = "82÷7=11, 5"

# Cell(9,1): "19÷9=2, 1" -> "71÷4=17, 3"
$tbl.Cell(9, 1).Range.Text = "71÷4=17, 3"

# Cell(9,2): "47÷5=9, 2" -> "28÷8=3, 4"
$tbl.Cell(9, 2).Range.Text = "28÷8=3, 4"

# Cell(9,3): "32÷3=10, 2" -> "67÷7=9, 4"
$tbl.Cell(9, 3).Range.Text = "67÷7=9, 4"

# Cell(9,4): "22÷3=7, 1" -> "82÷7=11, 5"
$tbl.Cell(9, 4).Range.Text = "82÷7=11, 5"

# Cell(9,5): "49÷2=24, 1" -> "50÷4=12, 2"
$tbl.Cell(9, 5).Range.Text = "50÷4=12, 2"

# Cell(13,1): "78÷7=11, 1" -> "62÷2=31, 0"
$tbl.Cell(13, 1).Range.Text = "62÷2=31, 0"

# Cell(13,2): "47÷9=5, 2" -> "43÷7=6, 1"
$tbl.Cell(13, 2).Range.Text = "43÷7=6, 1"

# Cell(13,3): "49÷7=7, 0" -> "17÷2=8, 1"
$tbl.Cell(13, 3).Range.Text = "17÷2=8, 1"

# Cell(13,4): "27÷2=13, 1" -> "23÷6=3, 5"
$tbl.Cell(13, 4).Range.Text = "23÷6=3, 5"

# Cell(13,5): "18÷3=6, 0" -> "74÷6=12, 2"
$tbl.Cell(13, 5).Range.Text = "74÷6=12, 2"

# Cell(17,1): "21÷8=2, 5" -> "49÷8=6, 1"
$tbl.Cell(17, 1).Range.Text = "49÷8=6, 1"

# Cell(17,2): "70÷3=23, 1" -> "10÷3=3, 1"
$tbl.Cell(17, 2).Range.Text = "10÷3=3, 1"

# Cell(17,3): "60÷9=6, 6" -> "96÷4=24, 0"
$tbl.Cell(17, 3).Range.Text = "96÷4=24, 0"

# Cell(17,4): "73÷8=9, 1" -> "66÷7=9, 3"
$tbl.Cell(17, 4).Range.Text = "66÷7=9, 3"

# Cell(17,5): "48÷3=16, 0" -> "87÷3=29, 0"
$tbl.Cell(17, 5).Range.Text = "87÷3=29, 0"

Write-Output "Updated date and 25 table cells."
